$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date-format tweak on the "Date Completed" column ---------------------
# Rows 2-6 and 10 already hold dates; switch their display format from the
# old "mm-dd-yy" (numFmtId 14) to "d-mmm-yy" (numFmtId 15).
$ws.Range("D2:D6").NumberFormat = "d-mmm-yy"
$ws.Range("D10").NumberFormat = "d-mmm-yy"

# --- Lesson 9 (row 10) recording link --------------------------------------
# Copy the "Recording" column's hyperlink-ish formatting from a row that
# already has it, then fill in the new Zoom recording link + value.
$ws.Range("E2").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newRecordingUrl = "https://zoom.us/recording/share/0cB7bdNC22KwtqFfrSVZ00CHZ5x5LZZ-Yna67qWQYMqwIumekTziMw?startTime=1561741589000"
$ws.Range("E10").Value = $newRecordingUrl
$ws.Hyperlinks.Add($ws.Range("E10"), $newRecordingUrl)

# --- Docker lesson (row 18) marked completed -------------------------------
$ws.Range("C18").Value = "Completed"
$ws.Range("D18").Value = 43648
$ws.Range("D18").NumberFormat = "d-mmm-yy"

# --- Update the on-screen selection ----------------------------------------
$ws.Range("B18").Select()
